$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.999.07"
$ws.Range("E2").Value = "  +0.34%  "
$ws.Range("D3").Value = "2.090.18"
$ws.Range("E3").Value = "  +2.93%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "'228.64"
$ws.Range("E5").Value = "  +0.24%  "
$ws.Range("E6").Value = "  +0.22%  "
$ws.Range("D7").Value = "'61.05"
$ws.Range("E7").Value = "  +0.96%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").Value = "'0.380"
$ws.Range("E9").Value = "  +0.07%  "
$ws.Range("D10").Value = "'0.0840"
$ws.Range("E10").Value = "  +2.75%  "
$ws.Range("E11").Value = "  +0.11%  "
$ws.Range("D12").Value = "2.398.67"
$ws.Range("E12").Value = "  +2.85%  "
$ws.Range("D13").Value = "'14.63"
$ws.Range("E13").Value = "  +0.54%  "
$ws.Range("D14").Value = "'22.13"
$ws.Range("E14").Value = "  +3.37%  "
$ws.Range("D15").Value = "'5.49"
$ws.Range("E15").Value = "  +6.33%  "
$ws.Range("E16").Value = "  +1.84%  "
$ws.Range("D17").Value = "2.065.91"
$ws.Range("E17").Value = "  +2.43%  "
$ws.Range("D18").Value = "37.609.88"
$ws.Range("E18").Value = "  -0.58%  "
$ws.Range("E19").Value = "  +1.97%  "
$ws.Range("D20").Value = "'69.98"
$ws.Range("E20").Value = "  +0.19%  "
$ws.Range("D21").Value = "0.0₃0839"
$ws.Range("E21").Value = "  +1.51%  "
$ws.Range("D22").Value = "'224.19"
$ws.Range("E22").Value = "  -0.01%  "
$ws.Range("E23").Value = "  +0.52%  "
$ws.Range("E24").Value = "  +0.09%  "
$ws.Range("E25").Value = "  +2.83%  "
$ws.Range("D26").Value = "'169.71"
$ws.Range("E26").Value = "  +1.65%  "
$ws.Range("D27").Value = "'9.37"
$ws.Range("E27").Value = "  +0.76%  "
$ws.Range("E28").Value = "  +3.30%  "
$ws.Range("D29").Value = "'18.96"
$ws.Range("E29").Value = "  +0.39%  "
$ws.Range("D30").Value = "'1.33"
$ws.Range("E30").Value = "  +4.09%  "
$ws.Range("E31").Value = "  -0.21%  "
$ws.Range("E32").Value = "  +10.58%  "
$ws.Range("B33").Value = "InternetComputer(DFINITY)"
$ws.Range("C33").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D33").Value = "'4.66"
$ws.Range("E33").Value = "  +3.64%  "
$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").Value = "'4.43"
$ws.Range("E34").Value = "  +0.50%  "
$ws.Range("E36").Value = "  +6.09%  "
$ws.Range("D37").Value = "'6.46"
$ws.Range("E37").Value = "  -0.13%  "
$ws.Range("E38").Value = "  +8.70%  "
$ws.Range("D39").Value = "'1.00"
$ws.Range("E39").Value = "  +0.07%  "
$ws.Range("D40").Value = "'18.02"
$ws.Range("E40").Value = "  +4.82%  "
$ws.Range("D41").Value = "1.545.63"
$ws.Range("E41").Value = "  +1.34%  "
$ws.Range("D42").Value = "'99.90"
$ws.Range("E42").Value = "  +4.15%  "
$ws.Range("D43").Value = "'0.0218"
$ws.Range("E43").Value = "  +0.45%  "
$ws.Range("E44").Value = "  -0.26%  "
$ws.Range("E45").Value = "  -0.78%  "
$ws.Range("D46").Value = "'4.16"
$ws.Range("E46").Value = "  +4.12%  "
$ws.Range("E47").Value = "  +0.74%  "
$ws.Range("E48").Value = "  +1.30%  "
$ws.Range("D49").Value = "'7.24"
$ws.Range("E49").Value = "  +1.85%  "
$ws.Range("E50").Value = "  +1.41%  "
$ws.Range("D51").Value = "2.287.27"
$ws.Range("E51").Value = "  +2.96%  "
